$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C5").Value = "WUKY"
$ws.Range("C6").Value = "WUKY"
$ws.Range("E11").Select() | Out-Null
